$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 113-114 (shifts existing rows 113-138 down to 115-140)
$rows = $ws.Range("A113:A114").EntireRow
$rows.Insert()

# New row 113: Primera, volume 60
$ws.Cells.Item(113, 1).Value = 7
$ws.Cells.Item(113, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(113, 3).Value = "Ñuble"
$ws.Cells.Item(113, 4).Value = "9/24/2021"
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100108
$ws.Cells.Item(113, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(113, 9).Value = 100108005
$ws.Cells.Item(113, 10).Value = "Piña"
$ws.Cells.Item(113, 11).Value = "Caramelo"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 60
$ws.Cells.Item(113, 14).Value = 19000
$ws.Cells.Item(113, 15).Value = 20000
$ws.Cells.Item(113, 16).Value = 19500
$ws.Cells.Item(113, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(113, 18).Value = "Ecuador"
$ws.Cells.Item(113, 19).Value = 1625
$ws.Cells.Item(113, 20).Value = 12

# New row 114: Segunda, volume 60
$ws.Cells.Item(114, 1).Value = 7
$ws.Cells.Item(114, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(114, 3).Value = "Ñuble"
$ws.Cells.Item(114, 4).Value = "9/24/2021"
$ws.Cells.Item(114, 5).Value = 16
$ws.Cells.Item(114, 6).Value = "Fruta"
$ws.Cells.Item(114, 7).Value = 100108
$ws.Cells.Item(114, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(114, 9).Value = 100108005
$ws.Cells.Item(114, 10).Value = "Piña"
$ws.Cells.Item(114, 11).Value = "Caramelo"
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 60
$ws.Cells.Item(114, 14).Value = 19000
$ws.Cells.Item(114, 15).Value = 20000
$ws.Cells.Item(114, 16).Value = 19500
$ws.Cells.Item(114, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(114, 18).Value = "Ecuador"
$ws.Cells.Item(114, 19).Value = 1393
$ws.Cells.Item(114, 20).Value = 14

Write-Host "Done"
